$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 and 5 (the last two data rows are removed entirely)
$ws.Range("A4:A5").EntireRow.Delete()

# Row 2 updates (corrected/recalculated TPM values)
$ws.Range("H2").Value = 0.651114
$ws.Range("I2").Value = 0.8572432933444277
$ws.Range("J2").Value = 0.8572432933444277
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.046121009076
$ws.Range("S2").Value = 0.8572432933444277
$ws.Range("T2").Value = 0.8572432933444277

# Row 3: swap Sending/Target cluster labels (FAPs <-> MuSCs) and update values
$ws.Range("A3").Value = "MuSCs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 0.03614333333333333
$ws.Range("H3").Value = 0.10843
$ws.Range("I3").Value = 0.1427567066555723
$ws.Range("J3").Value = 0.1427567066555723
$ws.Range("M3").Value = 0.212502
$ws.Range("N3").Value = 0.637506
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.00768053062
$ws.Range("R3").Value = 0.06912477558000001
$ws.Range("S3").Value = 0.1427567066555723
$ws.Range("T3").Value = 0.1427567066555723
